{"js": "// P2P connection works :)\n//\n// Rewrites the \"If user2 accepts ...\" bullet so that it documents the\n// full accept handshake: user2 sends (acc,user1uid) to the server, the\n// server relays (acc,user2uid) to user1, and only then does the server\n// connect the two peers. Also switches \"acc,\" -> \"inf,\" and \":userport\"\n// -> \",userport\" in the final \"sending each user (...)\" tuple, and moves\n// the trailing \"_GoBack\" bookmark to sit right after the rewritten tuple\n// (matching where Word's cursor last was).\n\nconst body = context.document.body;\n\n// 1) \" the server connects betwee\" -> the long expanded handshake text,\n//    ending in the same \"connects betwee\" tail so the following run\n//    (\"n them by\") still reads naturally.\nlet target1 = body.search(\" the server connects betwee\", { matchCase: true, ignorePunct: false });\ntarget1.load(\"items\");\nawait context.sync();\nif (target1.items.length > 0) {\n  target1.items[0].insertText(\n    \" (acc,user1uid) to the server. The server sends (acc,user2uid) to user1 and connects betwee\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) \"sending each user (acc,\" -> \"sending each user (inf,\"\nlet target2 = body.search(\"sending each user (acc,\", { matchCase: true, ignorePunct: false });\ntarget2.load(\"items\");\nawait context.sync();\nif (target2.items.length > 0) {\n  target2.items[0].insertText(\"sending each user (inf,\", \"Replace\");\n  await context.sync();\n}\n\n// 3) \":userport\" -> \",userport\"\nlet target3 = body.search(\":userport\", { matchCase: true, ignorePunct: false });\ntarget3.load(\"items\");\nawait context.sync();\nif (target3.items.length > 0) {\n  target3.items[0].insertText(\",userport\", \"Replace\");\n  await context.sync();\n}\n\n// 4) Move the \"_GoBack\" bookmark from the end of the \"If user2 rejects...\"\n//    paragraph to right after \"...userport\" (before the closing \")\") in\n//    the paragraph we just edited.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nlet target4 = body.search(\",userport\", { matchCase: true, ignorePunct: false });\ntarget4.load(\"items\");\nawait context.sync();\nif (target4.items.length > 0) {\n  const afterPort = target4.items[0].getRange(\"End\");\n  afterPort.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# P2P connection works :)\n#\n# Rewrites the \"If user2 accepts ...\" bullet so that it documents the\n# full accept handshake: user2 sends (acc,user1uid) to the server, the\n# server relays (acc,user2uid) to user1, and only then does the server\n# connect the two peers. Also switches \"acc,\" -> \"inf,\" and \":userport\"\n# -> \",userport\" in the final \"sending each user (...)\" tuple, and moves\n# the trailing \"_GoBack\" bookmark to sit right after the rewritten tuple\n# (matching where Word's cursor last was).\n\n$d = $word.ActiveDocument\n\n# 1) \" the server connects betwee\" -> expanded handshake text, still\n#    ending in \"connects betwee\" so the next run (\"n them by\") reads on.\n$r1 = $d.Content\n$r1.Find.Execute(\n    \" the server connects betwee\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \" (acc,user1uid) to the server. The server sends (acc,user2uid) to user1 and connects betwee\",\n    2\n)\n\n# 2) \"sending each user (acc,\" -> \"sending each user (inf,\"\n$r2 = $d.Content\n$r2.Find.Execute(\n    \"sending each user (acc,\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"sending each user (inf,\",\n    2\n)\n\n# 3) \":userport\" -> \",userport\"\n$r3 = $d.Content\n$r3.Find.Execute(\n    \":userport\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \",userport\",\n    2\n)\n\n# 4) Move the \"_GoBack\" bookmark from the end of the \"If user2 rejects...\"\n#    paragraph to right after \"...userport\" (before the closing \")\") in\n#    the paragraph we just edited.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Delete()\n\n$r4 = $d.Content\n$r4.Find.Execute(\n    \",userport\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"\",\n    0\n)\n$insertPoint = $d.Range($r4.End, $r4.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n"}
